# gg_fare_filing.xlsx edit script
# "use excel input to control which rbd is round trip only"
#
# 1. cabin_mapping: add a new "rt_only" table column, flagged Y for the
#    booking classes that are round-trip only (Y, B, M, K, R).
# 2. fare_combination: reorder columns so oneway_multiplier sits next to
#    weekend (B), weekend_surcharge moves to C, oneway moves to D; the
#    oneway column picks up an explicit "O" marker for the one-way rows.
# 3. L / K1 / H1 / P: fare_basis codes drop the old "...SE" country suffix
#    in favor of a "...US"/"...WOUS" suffix reflecting the new RBD scheme.
# 4. Selection/active-tab bookkeeping to match what Excel leaves behind.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. cabin_mapping: add rt_only column
# ---------------------------------------------------------------------
$wsCabin = $wb.Worksheets.Item("cabin_mapping")
$loCabin = $wsCabin.ListObjects.Item("Table2")
$loCabin.ListColumns.Add() | Out-Null

$hdrCell = $wsCabin.Range("C1")
$hdrCell.Value2 = "rt_only"
$hdrCell.HorizontalAlignment = -4108
$hdrCell.VerticalAlignment = -4108

$cabinBody = $loCabin.ListColumns.Item(3).DataBodyRange
$cabinBody.HorizontalAlignment = -4108
$cabinBody.VerticalAlignment = -4108

$rtOnlyClasses = @("Y", "B", "M", "K", "R")
for ($r = 2; $r -le 18; $r++) {
    $bookingClass = $wsCabin.Cells.Item($r, 1).Value2
    if ($rtOnlyClasses -contains $bookingClass) {
        $wsCabin.Cells.Item($r, 3).Value2 = "Y"
    }
}

# ---------------------------------------------------------------------
# 2. fare_combination: move oneway_multiplier next to weekend
# ---------------------------------------------------------------------
$wsFare = $wb.Worksheets.Item("fare_combination")

$wsFare.Columns("D").Cut() | Out-Null
$wsFare.Columns("B").Insert() | Out-Null

$wsFare.Range("D4").Value2 = "O"
$wsFare.Range("D5").Value2 = "O"

# Touch the header row so the table's column metadata resyncs to the
# new physical order (names/ids follow the header cell text).
$wsFare.Range("B1").Value2 = $wsFare.Range("B1").Value2
$wsFare.Range("C1").Value2 = $wsFare.Range("C1").Value2
$wsFare.Range("D1").Value2 = $wsFare.Range("D1").Value2

# ---------------------------------------------------------------------
# 3. L / K1 / H1 / P: fare_basis code updates
# ---------------------------------------------------------------------
$fareBasisMaps = @{
    "L" = @{
        "JLXSE" = "JLXUS"; "JLWSE" = "JLWOUS";
        "CLXSE" = "CLXUS"; "CLWSE" = "CLWOUS";
        "DLXSE" = "DLXUS"; "DLWSE" = "DLWOUS";
    }
    "K1" = @{
        "JK1XSE" = "JKXUS"; "JK1WSE" = "JKWOUS";
        "CK1XSE" = "CKXUS"; "CK1WSE" = "CKWOUS";
        "DK1XSE" = "DKXUS"; "DK1WSE" = "DKWOUS";
    }
    "H1" = @{
        "JH1XSE" = "JHXUS"; "JH1WSE" = "JHWOUS";
        "CH1XSE" = "CHXUS"; "CH1WSE" = "CHWOUS";
        "DH1XSE" = "DHXUS"; "DH1WSE" = "DHWOUS";
    }
    "P" = @{
        "JPXSE" = "JPXUS"; "JPWSE" = "JPWOUS";
        "CPXSE" = "CPXUS"; "CPWSE" = "CPWOUS";
        "DPXSE" = "DPXUS"; "DPWSE" = "DPWOUS";
    }
}

foreach ($sheetName in $fareBasisMaps.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $map = $fareBasisMaps[$sheetName]
    for ($r = 2; $r -le 13; $r++) {
        $cell = $ws.Cells.Item($r, 3)
        $old = $cell.Value2
        if ($map.ContainsKey($old)) {
            $cell.Value2 = $map[$old]
        }
    }
}

# ---------------------------------------------------------------------
# 4. Selection / active sheet bookkeeping
# ---------------------------------------------------------------------
$wsFare.Activate()
$wsFare.Range("C11").Select()

$wsCabin.Activate()
$wsCabin.Range("C16").Select()
